$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) values increased for a handful of events.
# These updates apply identically to both the "展览" and "全部类型" sheets,
# which mirror the same underlying data.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 2886
    $ws.Range("F4").Value = 99
    $ws.Range("F5").Value = 6699
    $ws.Range("F6").Value = 1615
    $ws.Range("F10").Value = 105
}
